$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Insert a new "Generic" task row right after the FEST-Logging "case
# study" row (old row 4), then fill in the task text before the section
# name (matches the order the strings were originally typed in) ---
$ws.Rows.Item(5).Insert()
$ws.Range("B5").Value = "change the UMLs to Visio PDFs"

# --- Mark a batch of tasks as DONE ---
$ws.Range("C2").Value = "DONE"
$ws.Range("C3").Value = "DONE"
$ws.Range("C13").Value = "DONE"
$ws.Range("C16").Value = "DONE"

# --- finish the new "Generic" row ---
$ws.Range("A5").Value = "Generic"

# --- Insert three more new rows at the top of the backlog (each insert
# pushes the previous new row down one slot) ---
$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "Drag-and-drop"
$ws.Range("B2").Value = "explain how sim. drag and drop works"

$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "Overview"
$ws.Range("B2").Value = "rework the overview"

$ws.Rows.Item(2).Insert()
$ws.Range("A2").Value = "Abstract"
$ws.Range("B2").Value = "rework the abstract"

# --- Table now spans 19 rows; drop the autofilter/sort leftovers and
# resize the table definition to match ---
$lo.ShowAutoFilter = $false
$lo.Resize($ws.Range("A1:C19"))

# --- Leave the selection where the author's cursor ended up ---
$ws.Range("B10").Select()
